# Generate Report for Archive
# - Update localization status text from "Ready for handoff" to "In Translation"
#   on the Overview sheet (zh-cn / de-de status columns) and on each language
#   sheet's Status column.
# - Narrow the "Latest HO Xliff Generate Date" status columns that used to be
#   sized for the old, longer status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update the status values that used to read "Ready for handoff".
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Shrink the now-narrower status columns to fit the shorter text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
